$d = $word.ActiveDocument

# 1. "Contenu attrayant" paragraph body text
$d.Content.Find.Execute(" le reel viral a réussi à transmettre l’appel du produit via des visuels attrayants et du contenu informatif.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " le reel viral a réussi à mettre en valeur les atouts du produit via des visuels attrayants et du contenu informatif.", 2)

# 2. "Marketing d’influenceurs" -> "Marketing d’influence" label
$d.Content.Find.Execute("Marketing d’influenceurs :", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Marketing d’influence :", 2)

# 2b. Marketing paragraph body text
$d.Content.Find.Execute(" le pouvoir du marketing d’influenceurs ne peut pas être surévalué.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " le pouvoir du marketing d’influence ne saurait être sous-estimé.", 2)

# 3. "Goût et variétés de saveurs" -> "Saveur et variété" label
$d.Content.Find.Execute("Goût et variétés de saveurs :", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Saveur et variété :", 2)

# 3b. Flavor paragraph body text
$d.Content.Find.Execute(" la réputation de Contoso Protein Plus pour les saveurs délicieuses et diverses était un point de vente clé dans le contenu viral.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " les saveurs délicieuses et variées qui font la réputation de Contoso Protein Plus ont été un argument de vente déterminant dans le contenu viral.", 2)

# 4. "Tendances en matière de santé et de fitness" -> "Tendances en matière de santé et de remise en forme" label
$d.Content.Find.Execute("Tendances en matière de santé et de fitness :", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tendances en matière de santé et de remise en forme :", 2)

# 4b. Health/fitness paragraph body text
$d.Content.Find.Execute(" la hausse continue de la conscience de la santé et de la forme physique, combinée à une augmentation du nombre de personnes qui adoptent des routines d’entraînement et des modes de vie actifs, a créé un marché réceptive pour un produit comme Contoso Protein Plus.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " l’intérêt croissant pour une vie saine et active, ainsi qu’une augmentation du nombre de personnes qui suivent des programmes de fitness et adoptent des modes de vie actifs, ont créé un marché réceptif à un produit comme Contoso Protein Plus.", 2)

# 5. "Disponibilité facile" -> "Disponibilité" label
$d.Content.Find.Execute("Disponibilité facile :", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Disponibilité :", 2)

# 5b. Availability paragraph body text
$d.Content.Find.Execute(" l’accessibilité du produit par le biais de différents détaillants en ligne a davantage alimenté la hype.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " la disponibilité du produit chez différents détaillants en ligne explique aussi son succès.", 2)

# 6. "Commentaires positifs et témoignages" -> "Commentaires et témoignages positifs" label
$d.Content.Find.Execute("Commentaires positifs et témoignages" + [char]0xa0 + ":", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Commentaires et témoignages positifs :", 2)

# 6b. Testimonials paragraph body text
$d.Content.Find.Execute(" Le reel n’était pas un cas isolé.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " le reel n’était pas un cas isolé.", 2)

# 7. "Word of Mouth" -> "Bouche à oreille" label
$d.Content.Find.Execute("Word of Mouth :", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Bouche à oreille :", 2)

# 7b. Word of mouth paragraph body text
$d.Content.Find.Execute(" les plateformes de médias sociaux favorisent la propagation rapide des tendances par le biais de la bouche à oreille.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " les plateformes de médias sociaux favorisent la propagation rapide des tendances par le bouche à oreille.", 2)
